# Update column F ("dSF") values for the rows that changed in this data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    3  = -4
    4  = -4
    6  = -8
    7  = -1
    10 = -1
    11 = -1
    12 = -7
    13 = -1
    14 = 4
    17 = -2
    20 = 1
    21 = -2
    22 = -2
    25 = 5
    26 = -2
    28 = -3
    29 = -4
    35 = 0
    36 = 5
    39 = -3
    40 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
